# Auto-generated script to apply Belias_Profits.xlsx market-data refresh
# Updates columns H-N (market price / leve profit calcs) for specific rows
# across sheets ALC, ARM, BSM, CRP, GSM, LTW, per the scheduled runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2033.25
$ws.Range("I40").Value = 1066.6666
$ws.Range("J40").Value = 2355.4443
$ws.Range("K40").Value = 1066.6666
$ws.Range("L40").Value = 2355.4443
$ws.Range("M40").Value = -891.6666
$ws.Range("N40").Value = -2705.4443
# Row 64
$ws.Range("H64").Value = 3073
$ws.Range("I64").Value = 2699
$ws.Range("K64").Value = 2699
$ws.Range("M64").Value = -2451
# Row 67
$ws.Range("H67").Value = 3073
$ws.Range("I67").Value = 2699
$ws.Range("K67").Value = 2699
$ws.Range("M67").Value = -1841
# Row 76
$ws.Range("H76").Value = 3307.5
$ws.Range("I76").Value = 3335.4546
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3335.4546
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -3020.4546
$ws.Range("N76").Value = -3630
# Row 79
$ws.Range("H79").Value = 3307.5
$ws.Range("I79").Value = 3335.4546
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3335.4546
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2243.4546
$ws.Range("N79").Value = -5184
# Row 100
$ws.Range("H100").Value = 2889.5557
$ws.Range("I100").Value = 2866.6667
$ws.Range("J100").Value = 2935.3333
$ws.Range("K100").Value = 2866.6667
$ws.Range("L100").Value = 2935.3333
$ws.Range("M100").Value = -2325.6667
$ws.Range("N100").Value = -4017.3333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22220.94
$ws.Range("I32").Value = 9308.807
$ws.Range("J32").Value = 33339.723
$ws.Range("K32").Value = 9308.807
$ws.Range("L32").Value = 33339.723
$ws.Range("M32").Value = -9021.807
$ws.Range("N32").Value = -33913.723
# Row 63
$ws.Range("H63").Value = 3707.4167
$ws.Range("I63").Value = 2320.2
$ws.Range("J63").Value = 4698.2856
$ws.Range("K63").Value = 2320.2
$ws.Range("L63").Value = 4698.2856
$ws.Range("M63").Value = -1634.2
$ws.Range("N63").Value = -6070.2856
# Row 66
$ws.Range("H66").Value = 3707.4167
$ws.Range("I66").Value = 2320.2
$ws.Range("J66").Value = 4698.2856
$ws.Range("K66").Value = 11601
$ws.Range("L66").Value = 23491.428
$ws.Range("M66").Value = -8169
$ws.Range("N66").Value = -30355.428
# Row 132
$ws.Range("H132").Value = 1396.4147
$ws.Range("I132").Value = 1101.6666
$ws.Range("J132").Value = 2612.25
$ws.Range("K132").Value = 3304.9998
$ws.Range("L132").Value = 7836.75
$ws.Range("M132").Value = -774.9998
$ws.Range("N132").Value = -12896.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2175
$ws.Range("I105").Value = 1850
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1850
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -103
$ws.Range("N105").Value = -5994

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2247.5574
$ws.Range("I31").Value = 1016.8823
$ws.Range("J31").Value = 3797.2964
$ws.Range("K31").Value = 1016.8823
$ws.Range("L31").Value = 3797.2964
$ws.Range("M31").Value = -721.8823
$ws.Range("N31").Value = -4387.2964
# Row 34
$ws.Range("H34").Value = 2247.5574
$ws.Range("I34").Value = 1016.8823
$ws.Range("J34").Value = 3797.2964
$ws.Range("K34").Value = 1016.8823
$ws.Range("L34").Value = 3797.2964
$ws.Range("M34").Value = -814.8823
$ws.Range("N34").Value = -4201.2964
# Row 62
$ws.Range("H62").Value = 2481.25
$ws.Range("I62").Value = 2433.3333
$ws.Range("J62").Value = 2625
$ws.Range("K62").Value = 2433.3333
$ws.Range("L62").Value = 2625
$ws.Range("M62").Value = -1809.3333
$ws.Range("N62").Value = -3873
# Row 65
$ws.Range("H65").Value = 2481.25
$ws.Range("I65").Value = 2433.3333
$ws.Range("J65").Value = 2625
$ws.Range("K65").Value = 12166.6665
$ws.Range("L65").Value = 13125
$ws.Range("M65").Value = -9046.6665
$ws.Range("N65").Value = -19365
# Row 134
$ws.Range("H134").Value = 2665.4707
$ws.Range("I134").Value = 2578.318
$ws.Range("J134").Value = 2825.25
$ws.Range("K134").Value = 7734.954000000001
$ws.Range("L134").Value = 8475.75
$ws.Range("M134").Value = -5199.954000000001
$ws.Range("N134").Value = -13545.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
# Row 70
$ws.Range("H70").Value = 5298.1665
$ws.Range("I70").Value = 4717.3335
$ws.Range("J70").Value = 6459.8335
$ws.Range("K70").Value = 4717.3335
$ws.Range("L70").Value = 6459.8335
$ws.Range("M70").Value = -4447.3335
$ws.Range("N70").Value = -6999.8335
# Row 73
$ws.Range("H73").Value = 5298.1665
$ws.Range("I73").Value = 4717.3335
$ws.Range("J73").Value = 6459.8335
$ws.Range("K73").Value = 4717.3335
$ws.Range("L73").Value = 6459.8335
$ws.Range("M73").Value = -3781.3335
$ws.Range("N73").Value = -8331.8335
# Row 80
$ws.Range("H80").Value = 2414.2856
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -4496
# Row 83
$ws.Range("H83").Value = 2414.2856
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -22484
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
# Row 112
$ws.Range("H112").Value = 28000
$ws.Range("J112").Value = 28000
$ws.Range("L112").Value = 28000
$ws.Range("N112").Value = -30216
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 5499.2607
$ws.Range("I132").Value = 5251.722
$ws.Range("J132").Value = 6390.4
$ws.Range("K132").Value = 15755.166
$ws.Range("L132").Value = 19171.2
$ws.Range("M132").Value = -13225.166
$ws.Range("N132").Value = -24231.2
